$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text/trait edits on existing rows ----
$ws.Range("Y29").Value = 'Children Servants, Divine, Greek Mythology Males, Humanoid, King, Male, Riding, Servant, Weak to Enuma Elish'
$ws.Range("Y34").Value = 'Children Servants, Humanoid, Male, Servant, Weak to Enuma Elish'
$ws.Range("Y54").Value = 'Children Servants, Giant, Greek Mythology Males, Humanoid, Male, Servant, Weak to Enuma Elish'
$ws.Range("Y57").Value = 'Children Servants, Dragon, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y75").Value = 'Children Servants (Stage2-3), Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y76").Value = 'Children Servants, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y96").Value = 'Children Servants, Divine, Humanoid, King, Male, Servant, Weak to Enuma Elish'
$ws.Range("Y118").Value = 'Children Servants, Demonic, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y136").Value = 'Children Servants, Dragon, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y138").Value = 'Children Servants, Female, Humanoid, Illya, Living Human, Pseudo-Servants, Servant, Weak to Enuma Elish'
$ws.Range("Y139").Value = 'Children Servants, Female, Humanoid, Living Human, Servant, Weak to Enuma Elish'
$ws.Range("Y143").Value = 'Children Servants, Female, Humanoid, Saberface, Servant, Weak to Enuma Elish'
$ws.Range("Y148").Value = 'Children Servants, Divine, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y161").Value = 'Children Servants, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y168").Value = 'Children Servants, Female, Humanoid, King, Servant, Weak to Enuma Elish'
$ws.Range("Y172").Value = 'Children Servants, Female, Giant, Humanoid, Servants, Weak to Enuma Elish'
$ws.Range("Y193").Value = 'Children Servants, Divine, Female, Humanoid, Servant, Threat to Humanity'
$ws.Range("F208").Value = 'Manslayer Izō (人斬り以蔵, Hitokiri Izō?), Okada Izou, Ghost of Tosa'
$ws.Range("Y215").Value = 'Children Servants, Demonic, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y222").Value = 'Children Servants, Divine, Female, Humanoid, Illya, Pseudo-Servant, Servant, Weak to Enuma Elish'
$ws.Range("A228").Value = 'Consort Yu'
$ws.Range("F228").Value = 'Yu Miaoyi, Yu the Beauty, Yu Meiren, Akuta Hinako (芥ヒナコ?)Crimson Beauty Under the Moon (紅の月下美人?), Gubijin, Gucchan, 虞妙弋'
$ws.Range("Y232").Value = 'Children Servants, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Y234").Value = 'Children Servants, Female, Humanoid, Living Human, Pseudo-Servants, Servant, Weak to Enuma Elish'
$ws.Range("Y237").Value = 'Children Servants (Stage 1), Divine, Female, Humanoid, Pseudo-Servant, Riding, Servant, Weak to Enuma Elish'
$ws.Range("Y246").Value = 'Argo-Related, Divine, Greek Mythology Males, Humanoid, Male, Servant, Weak to Enuma Elish'
$ws.Range("F250").Value = 'Da Vinci Lily, Loli Vinci, Gran Cavallo'
$ws.Range("Y250").Value = 'Children Servants, Female, Humanoid, Riding, Servant'
$ws.Range("Y252").Value = 'Argo-Related, Brynhildr''s Beloved, Children Servants, Greek Mythology Males, Humanoid, Male, Servant, Weak to Enuma Elish'
$ws.Range("AH269").Value = '－'
$ws.Range("Y278").Value = 'Children Servants, Humanoid, Male, Servant'
$ws.Range("Y280").Value = 'Children Servants, Divine, Female, Humanoid, Living Human, Riding, Servant, Weak to Enuma Elish'
$ws.Range("F281").Value = 'Castoria, CasSeiba'
$ws.Range("Y281").Value = 'Arthur, Humanoid, Saberface, Servant, Female'

# ---- New rows 282-287 ----
# Columns S:V hold percentage-looking text (e.g. "14.8%"); Excel would
# otherwise auto-coerce these into numeric percent cells on assignment,
# so force them to literal Text format first to match the source data.
$ws.Range("S282:V287").NumberFormat = "@"

# Row 282
$ws.Range("A282").Value = 'Sesshōin Kiara (Moon Cancer)'
$ws.Range("B282").Value = 285
$ws.Range("C282").Value = '5-Star'
$ws.Range("D282").Value = 'Moon Cancer'
$ws.Range("E282").Value = '殺生院キアラ'
$ws.Range("F282").Value = 'Sessyoin Kiara, Sesshouin Kiara, Demonic Bodhisattva, Last Prophet, Beast III/R'
$ws.Range("G282").Value = 16
$ws.Range("H282").Value = 1719
$ws.Range("I282").Value = 11128
$ws.Range("J282").Value = 2249
$ws.Range("K282").Value = 15336
$ws.Range("L282").Value = 12181
$ws.Range("M282").Value = 16801
$ws.Range("N282").Value = 'Tanaka Rie'
$ws.Range("O282").Value = 'Wada Arco'
$ws.Range("P282").Value = 'Earth'
$ws.Range("Q282").Value = 'S'
$ws.Range("R282").Value = 50
$ws.Range("S282").Value = '14.8%'
$ws.Range("T282").Value = '0.6%'
$ws.Range("U282").Value = '3%'
$ws.Range("V282").Value = '0.5%'
$ws.Range("W282").Value = 'Lawful Evil'
$ws.Range("X282").Value = 'Female'
$ws.Range("Y282").Value = 'Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Z282").Value = 'QAABB'
$ws.Range("AA282").Value = 3
$ws.Range("AB282").Value = 3
$ws.Range("AC282").Value = 4
$ws.Range("AD282").Value = 5
$ws.Range("AE282").Value = 'A'
$ws.Range("AF282").Value = 'Anti-World'
$ws.Range("AG282").Value = 'Arts'
$ws.Range("AH282").Value = 3

# Row 283
$ws.Range("A283").Value = 'Illyasviel von Einzbern (Archer)'
$ws.Range("B283").Value = 286
$ws.Range("C283").Value = '4-Star'
$ws.Range("D283").Value = 'Archer'
$ws.Range("E283").Value = 'イリヤスフィール・フォン・アインツベルン'
$ws.Range("F283").Value = 'Illya'
$ws.Range("G283").Value = 12
$ws.Range("H283").Value = 1683
$ws.Range("I283").Value = 10098
$ws.Range("J283").Value = 1746
$ws.Range("K283").Value = 10914
$ws.Range("L283").Value = 12226
$ws.Range("M283").Value = 13233
$ws.Range("N283").Value = 'Kadowaki Mai, Takano Naoko'
$ws.Range("O283").Value = 'Hiroyama Hiroshi'
$ws.Range("P283").Value = 'Man'
$ws.Range("Q283").Value = 'Semi Reverse S'
$ws.Range("R283").Value = 148
$ws.Range("S283").Value = '8%'
$ws.Range("T283").Value = '0.63%'
$ws.Range("U283").Value = '3%'
$ws.Range("V283").Value = '31.5%'
$ws.Range("W283").Value = 'Neutral Summer'
$ws.Range("X283").Value = 'Female'
$ws.Range("Y283").Value = 'Children Servants, Female, Humanoid, Illya, Living Human, Pseudo-Servants, Servant, Weak to Enuma Elish'
$ws.Range("Z283").Value = 'QQAAB'
$ws.Range("AA283").Value = 4
$ws.Range("AB283").Value = 3
$ws.Range("AC283").Value = 3
$ws.Range("AD283").Value = 4
$ws.Range("AE283").Value = 'A'
$ws.Range("AF283").Value = 'Anti-Army'
$ws.Range("AG283").Value = 'Quick'
$ws.Range("AH283").Value = 5

# Row 284
$ws.Range("A284").Value = 'Brynhildr (Berserker)'
$ws.Range("B284").Value = 287
$ws.Range("C284").Value = '4-Star'
$ws.Range("D284").Value = 'Berserker'
$ws.Range("E284").Value = 'ブリュンヒルデ'
$ws.Range("F284").Value = 'Brynhild (Brünnhilde), Sigrdrífa'
$ws.Range("G284").Value = 12
$ws.Range("H284").Value = 1699
$ws.Range("I284").Value = 10197
$ws.Range("J284").Value = 1603
$ws.Range("K284").Value = 10023
$ws.Range("L284").Value = 12346
$ws.Range("M284").Value = 12153
$ws.Range("N284").Value = 'Noto Mamiko'
$ws.Range("O284").Value = 'Miwa Shirō'
$ws.Range("P284").Value = 'Sky'
$ws.Range("Q284").Value = 'Semi Reverse S'
$ws.Range("R284").Value = 10
$ws.Range("S284").Value = '5.1%'
$ws.Range("T284").Value = '0.85%'
$ws.Range("U284").Value = '5%'
$ws.Range("V284").Value = '45.5%'
$ws.Range("W284").Value = 'Neutral Good'
$ws.Range("X284").Value = 'Female'
$ws.Range("Y284").Value = 'Divine, Female, Humanoid, Servant, Weak to Enuma Elish'
$ws.Range("Z284").Value = 'QABBB'
$ws.Range("AA284").Value = 3
$ws.Range("AB284").Value = 3
$ws.Range("AC284").Value = 4
$ws.Range("AD284").Value = 5
$ws.Range("AE284").Value = 'B'
$ws.Range("AF284").Value = 'Anti-Unit / Anti-Army'
$ws.Range("AG284").Value = 'Buster'
$ws.Range("AH284").Value = 5

# Row 285
$ws.Range("A285").Value = 'Consort Yu (Lancer)'
$ws.Range("B285").Value = 288
$ws.Range("C285").Value = '4-Star'
$ws.Range("D285").Value = 'Lancer'
$ws.Range("E285").Value = '虞美人'
$ws.Range("F285").Value = 'Yu Miaoyi, Yu the Beauty, Yu Meiren, Akuta Hinako (芥ヒナコ?)Crimson Beauty Under the Moon (紅の月下美人?), Gubijin, Gucchan, 虞妙弋'
$ws.Range("G285").Value = 12
$ws.Range("H285").Value = 1649
$ws.Range("I285").Value = 9896
$ws.Range("J285").Value = 1799
$ws.Range("K285").Value = 11245
$ws.Range("L285").Value = 11982
$ws.Range("M285").Value = 13634
$ws.Range("N285").Value = 'Ise Mariya'
$ws.Range("O285").Value = 'toi8'
$ws.Range("P285").Value = 'Earth'
$ws.Range("Q285").Value = 'Reverse S'
$ws.Range("R285").Value = 89
$ws.Range("S285").Value = '12.2%'
$ws.Range("T285").Value = '1.1%'
$ws.Range("U285").Value = '4%'
$ws.Range("V285").Value = '24%'
$ws.Range("W285").Value = 'Lawful Evil'
$ws.Range("X285").Value = 'Female'
$ws.Range("Y285").Value = 'Demonic, Humanoid, Female, Servant, Weak to Enuma Elish'
$ws.Range("Z285").Value = 'QQABB'
$ws.Range("AA285").Value = 3
$ws.Range("AB285").Value = 2
$ws.Range("AC285").Value = 3
$ws.Range("AD285").Value = 4
$ws.Range("AE285").Value = 'A'
$ws.Range("AF285").Value = 'Anti-Army'
$ws.Range("AG285").Value = 'Quick'
$ws.Range("AH285").Value = 5

# Row 286
$ws.Range("A286").Value = 'Abigail Williams (Summer)'
$ws.Range("B286").Value = 289
$ws.Range("C286").Value = '5-Star'
$ws.Range("D286").Value = 'Foreigner'
$ws.Range("E286").Value = 'アビゲイル・ウィリアムズ'
$ws.Range("F286").Value = 'The Key to the Gate, Sut-Typhon, Yog-Sothoth, All-in-One, One-in-All, Abby'
$ws.Range("G286").Value = 16
$ws.Range("H286").Value = 1820
$ws.Range("I286").Value = 11781
$ws.Range("J286").Value = 2090
$ws.Range("K286").Value = 14250
$ws.Range("L286").Value = 12896
$ws.Range("M286").Value = 15611
$ws.Range("N286").Value = 'Ōwada Hitomi'
$ws.Range("O286").Value = 'Kuroboshi Kouhaku'
$ws.Range("P286").Value = 'Earth'
$ws.Range("Q286").Value = 'Semi S'
$ws.Range("R286").Value = 150
$ws.Range("S286").Value = '15%'
$ws.Range("T286").Value = '0.86%'
$ws.Range("U286").Value = '3%'
$ws.Range("V286").Value = '6%'
$ws.Range("W286").Value = 'Lawful Evil'
$ws.Range("X286").Value = 'Female'
$ws.Range("Y286").Value = 'Children Servants, Divine, Female, Humanoid, Servant, Threat to Humanity'
$ws.Range("Z286").Value = 'QQABB'
$ws.Range("AA286").Value = 3
$ws.Range("AB286").Value = 3
$ws.Range("AC286").Value = 3
$ws.Range("AD286").Value = 5
$ws.Range("AE286").Value = 'A'
$ws.Range("AF286").Value = 'Anti-Unit'
$ws.Range("AG286").Value = 'Buster'
$ws.Range("AH286").Value = 3

# Row 287
$ws.Range("A287").Value = 'Tomoe Gozen (Saber)'
$ws.Range("B287").Value = 290
$ws.Range("C287").Value = '4-Star'
$ws.Range("D287").Value = 'Saber'
$ws.Range("E287").Value = '巴御前'
$ws.Range("F287").Value = '-'
$ws.Range("G287").Value = 12
$ws.Range("H287").Value = 1590
$ws.Range("I287").Value = 9544
$ws.Range("J287").Value = 1957
$ws.Range("K287").Value = 12233
$ws.Range("L287").Value = 11556
$ws.Range("M287").Value = 14832
$ws.Range("N287").Value = 'Kanemoto Hisako'
$ws.Range("O287").Value = 'Shirabi'
$ws.Range("P287").Value = 'Earth'
$ws.Range("Q287").Value = 'Linear'
$ws.Range("R287").Value = 102
$ws.Range("S287").Value = '10%'
$ws.Range("T287").Value = '0.55%'
$ws.Range("U287").Value = '3%'
$ws.Range("V287").Value = '24.5%'
$ws.Range("W287").Value = 'Neutral   Good'
$ws.Range("X287").Value = 'Female'
$ws.Range("Y287").Value = 'Demonic, Female, Humanoid, Riding, Servant, Weak to Enuma Elish'
$ws.Range("Z287").Value = 'QAABB'
$ws.Range("AA287").Value = 4
$ws.Range("AB287").Value = 3
$ws.Range("AC287").Value = 2
$ws.Range("AD287").Value = 4
$ws.Range("AE287").Value = 'B'
$ws.Range("AF287").Value = 'Anti-Army'
$ws.Range("AG287").Value = 'Arts'
$ws.Range("AH287").Value = 3

